$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36, shifting existing rows 36-44 down to 37-45
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new weekly record
$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(36, 3).Value = "Bíobío"
$ws.Cells.Item(36, 4).Value = 44474
$ws.Cells.Item(36, 5).Value = 8
$ws.Cells.Item(36, 6).Value = 100112012
$ws.Cells.Item(36, 7).Value = "Espinaca"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 50
$ws.Cells.Item(36, 11).Value = 6000
$ws.Cells.Item(36, 12).Value = 6500
$ws.Cells.Item(36, 13).Value = 6300
$ws.Cells.Item(36, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 630
$ws.Cells.Item(36, 17).Value = 10
$ws.Cells.Item(36, 18).Value = "Hortaliza"
